$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: ring-buffer summary numbers change ---
$ws.Range("B5").Value = 14077775
$ws.Range("C5").Value = 1075
$ws.Range("D5").Value = 2518
$ws.Range("E5").Value = 2235
$ws.Range("F5").Value = 2014
$ws.Range("G5").Value = 1364
$ws.Range("H5").Value = 1986
$ws.Range("I5").Value = 2618

# --- Rows 6:7 (old "a"/"b" duplicate rows) are removed entirely ---
$ws.Range("A6:I7").ClearContents()

# --- Row 11: (A) Locked -- unchanged values, kept for completeness ---
$ws.Range("A11").Value = "64 Core Intel Xeon CPU E7-4820 @ 2.00GHz (A) Locked"
$ws.Range("B11").Value = 6558774
$ws.Range("C11").Value = 8914003
$ws.Range("D11").Value = 8757513
$ws.Range("E11").Value = 8415938
$ws.Range("F11").Value = 7775145
$ws.Range("G11").Value = 4334278
$ws.Range("H11").Value = 3166165
$ws.Range("I11").Value = 3166369

# --- Row 12: (A) CAS lock ---
$ws.Range("A12").Value = "64 Core Intel Xeon CPU E7-4820 @ 2.00GHz (A) CAS lock"
$ws.Range("B12").Value = 5642544
$ws.Range("C12").Value = 8288701
$ws.Range("D12").Value = 5707577
$ws.Range("E12").Value = 5710434
$ws.Range("F12").Value = 7054031
$ws.Range("G12").Value = 5687563
$ws.Range("H12").Value = 7064434
$ws.Range("I12").Value = 5705941

# --- Row 13: (A) Ticket (new) ---
$ws.Range("A13").Value = "64 Core Intel Xeon CPU E7-4820 @ 2.00GHz (A) Ticket"
$ws.Range("B13").Value = 6163242
$ws.Range("C13").Value = 40341
$ws.Range("D13").Value = 15
$ws.Range("E13").Value = 29876
$ws.Range("F13").Value = 125211
$ws.Range("G13").Value = 39334
$ws.Range("H13").Value = 11527
$ws.Range("I13").Value = 39827

# --- Row 14: (A) TAS (new) ---
$ws.Range("A14").Value = "64 Core Intel Xeon CPU E7-4820 @ 2.00GHz (A) TAS"
$ws.Range("B14").Value = 4910844
$ws.Range("C14").Value = 5505544
$ws.Range("D14").Value = 2979812
$ws.Range("E14").Value = 715369
$ws.Range("F14").Value = 215518
$ws.Range("G14").Value = 130338
$ws.Range("H14").Value = 66332
$ws.Range("I14").Value = 36518

# --- Row 15: (A) TTAS (new, replaces old (A) C++ Spinlock row) ---
$ws.Range("A15").Value = "64 Core Intel Xeon CPU E7-4820 @ 2.00GHz (A) TTAS"
$ws.Range("B15").Value = 4235748
$ws.Range("C15").Value = 7986304
$ws.Range("D15").Value = 7533212
$ws.Range("E15").Value = 4266191
$ws.Range("F15").Value = 4239697
$ws.Range("G15").Value = 4910270
$ws.Range("H15").Value = 4221463
$ws.Range("I15").Value = 3673234

# --- Row 16: (B) Locked ---
$ws.Range("A16").Value = "4 Core Intel Core i5-2500K CPU @ 3.30GHz (B) Locked"
$ws.Range("B16").Value = 5328217
$ws.Range("C16").Value = 418741
$ws.Range("D16").Value = 473308
$ws.Range("E16").Value = 461109
$ws.Range("F16").Value = 448229
$ws.Range("G16").Value = 448645
$ws.Range("H16").Value = 436339
$ws.Range("I16").Value = 432057

# --- Row 17: (B) CAS lock ---
$ws.Range("A17").Value = "4 Core Intel Core i5-2500K CPU @ 3.30GHz (B) CAS lock"
$ws.Range("B17").Value = 9771618
$ws.Range("C17").Value = 9735489
$ws.Range("D17").Value = 9905043
$ws.Range("E17").Value = 9916571
$ws.Range("F17").Value = 9639092
$ws.Range("G17").Value = 9802904
$ws.Range("H17").Value = 9930743
$ws.Range("I17").Value = 9515278

# --- Row 18: (B) Ticket (new) ---
$ws.Range("A18").Value = "4 Core Intel Core i5-2500K CPU @ 3.30GHz (B) Ticket"
$ws.Range("B18").Value = 14077775
$ws.Range("C18").Value = 1075
$ws.Range("D18").Value = 2518
$ws.Range("E18").Value = 2235
$ws.Range("F18").Value = 2014
$ws.Range("G18").Value = 1364
$ws.Range("H18").Value = 1986
$ws.Range("I18").Value = 2618

# --- Row 19: (B) TAS (new) ---
$ws.Range("A19").Value = "4 Core Intel Core i5-2500K CPU @ 3.30GHz (B) TAS"
$ws.Range("B19").Value = 10644755
$ws.Range("C19").Value = 8292727
$ws.Range("D19").Value = 4840726
$ws.Range("E19").Value = 4499906
$ws.Range("F19").Value = 4600061
$ws.Range("G19").Value = 4661298
$ws.Range("H19").Value = 4934930
$ws.Range("I19").Value = 4789439

# --- Row 20: (B) TTAS (new, replaces old (B) C++ Spinlock row) ---
$ws.Range("A20").Value = "4 Core Intel Core i5-2500K CPU @ 3.30GHz (B) TTAS"
$ws.Range("B20").Value = 10604337
$ws.Range("C20").Value = 10611323
$ws.Range("D20").Value = 10564897
$ws.Range("E20").Value = 10540146
$ws.Range("F20").Value = 10570471
$ws.Range("G20").Value = 10592364
$ws.Range("H20").Value = 10539685
$ws.Range("I20").Value = 10488351

# --- Row 21: (C - Ducss) Locked ---
$ws.Range("A21").Value = "2 Core Intel CPU @ 3.00 GHz (C - Ducss) Locked"
$ws.Range("B21").Value = 1321542
$ws.Range("C21").Value = 1495327
$ws.Range("D21").Value = 1467956
$ws.Range("E21").Value = 1254631
$ws.Range("F21").Value = 1363862
$ws.Range("G21").Value = 1767556
$ws.Range("H21").Value = 1358405
$ws.Range("I21").Value = 1373794

# --- Row 22: (C) CAS lock ---
$ws.Range("A22").Value = "2 Core Intel CPU @ 3.00 GHz (C) CAS lock"
$ws.Range("B22").Value = 1249352
$ws.Range("C22").Value = 1258630
$ws.Range("D22").Value = 1438296
$ws.Range("E22").Value = 1247512
$ws.Range("F22").Value = 1438910
$ws.Range("G22").Value = 1230792
$ws.Range("H22").Value = 1438801
$ws.Range("I22").Value = 1233541

# --- Row 23: (C) Ticket (new row) ---
$ws.Range("A23").Value = "2 Core Intel CPU @ 3.00 GHz (C)  Ticket"
$ws.Range("B23").Value = 1401147
$ws.Range("C23").Value = 40
$ws.Range("D23").Value = 107
$ws.Range("E23").Value = 265
$ws.Range("F23").Value = 544
$ws.Range("G23").Value = 1139
$ws.Range("H23").Value = 2445
$ws.Range("I23").Value = 5137

# --- Row 24: (C) TAS (new row) ---
$ws.Range("A24").Value = "2 Core Intel CPU @ 3.00 GHz (C)  TAS"
$ws.Range("B24").Value = 1303082
$ws.Range("C24").Value = 1850386
$ws.Range("D24").Value = 1327536
$ws.Range("E24").Value = 824253
$ws.Range("F24").Value = 413418
$ws.Range("G24").Value = 245087
$ws.Range("H24").Value = 167994
$ws.Range("I24").Value = 51203

# --- Row 25: (C) TTAS (new row, replaces old (C) Spinlock row) ---
$ws.Range("A25").Value = "2 Core Intel CPU @ 3.00 GHz (C)  TTAS"
$ws.Range("B25").Value = 1191029
$ws.Range("C25").Value = 1180264
$ws.Range("D25").Value = 1123372
$ws.Range("E25").Value = 1189681
$ws.Range("F25").Value = 1200067
$ws.Range("G25").Value = 1228259
$ws.Range("H25").Value = 1205112
$ws.Range("I25").Value = 1179244

# --- Row 26: (D) Locked ---
$ws.Range("A26").Value = "2 Core Intel CPU @ 2.80 GHz (D) Locked"
$ws.Range("B26").Value = 3421422
$ws.Range("C26").Value = 4856188
$ws.Range("D26").Value = 4855373
$ws.Range("E26").Value = 4557720
$ws.Range("F26").Value = 4723686
$ws.Range("G26").Value = 3093845
$ws.Range("H26").Value = 3169190
$ws.Range("I26").Value = 2640398

# --- Row 27: (D) CAS lock ---
$ws.Range("A27").Value = "2 Core Intel CPU @ 2.80 GHz (D) CAS lock"
$ws.Range("B27").Value = 2642745
$ws.Range("C27").Value = 2716931
$ws.Range("D27").Value = 3473391
$ws.Range("E27").Value = 2658225
$ws.Range("F27").Value = 2654354
$ws.Range("G27").Value = 3450999
$ws.Range("H27").Value = 2577043
$ws.Range("I27").Value = 2596864

# --- Row 28: (D) Ticket (new row) ---
$ws.Range("A28").Value = "2 Core Intel CPU @ 2.80 GHz (D) Ticket"
$ws.Range("B28").Value = 3960868
$ws.Range("C28").Value = 852
$ws.Range("D28").Value = 370
$ws.Range("E28").Value = 684
$ws.Range("F28").Value = 1696
$ws.Range("G28").Value = 6064
$ws.Range("H28").Value = 142
$ws.Range("I28").Value = 16590

# --- Row 29: (D) TAS (new row) ---
$ws.Range("A29").Value = "2 Core Intel CPU @ 2.80 GHz (D) TAS"
$ws.Range("B29").Value = 2581570
$ws.Range("C29").Value = 2590346
$ws.Range("D29").Value = 2504325
$ws.Range("E29").Value = 1355264
$ws.Range("F29").Value = 616020
$ws.Range("G29").Value = 252520
$ws.Range("H29").Value = 74205
$ws.Range("I29").Value = 29459

# --- Row 30: (D) TTAS (new row, replaces old (D) Spinlock row) ---
$ws.Range("A30").Value = "2 Core Intel CPU @ 2.80 GHz (D) TTAS"
$ws.Range("B30").Value = 2560144
$ws.Range("C30").Value = 3529022
$ws.Range("D30").Value = 2567171
$ws.Range("E30").Value = 2586400
$ws.Range("F30").Value = 2599223
$ws.Range("G30").Value = 2575840
$ws.Range("H30").Value = 2522506
$ws.Range("I30").Value = 2445529

# --- Update selection / active cell to match the new view state ---
$ws.Range("A10:I30").Select()
